$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.023199999999997
$ws.Range("B9").Value = 8.628400000000006
$ws.Range("B18").Value = 4.851000000000003
$ws.Range("B20").Value = 5.672299999999995
$ws.Range("E21").Value = 13.283
